$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 121 - this pushes existing rows 121..155
# down to 122..156, preserving all of their data/formatting intact
# (matches the diff: old row 121 data now lives at row 122, ..., old row
# 155 data now lives at row 156).
$ws.Rows.Item(121).Insert()

# Populate the newly inserted row 121 with the new record.
$ws.Cells.Item(121, 1).Value = 11
$ws.Cells.Item(121, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(121, 3).Value = "Bíobío"
$ws.Cells.Item(121, 4).Value = 44642
$ws.Cells.Item(121, 5).Value = 8
$ws.Cells.Item(121, 6).Value = "Fruta"
$ws.Cells.Item(121, 7).Value = 100108
$ws.Cells.Item(121, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(121, 9).Value = 100108005
$ws.Cells.Item(121, 10).Value = "Piña"
$ws.Cells.Item(121, 11).Value = "Sin especificar"
$ws.Cells.Item(121, 12).Value = "Segunda"
$ws.Cells.Item(121, 13).Value = 270
$ws.Cells.Item(121, 14).Value = 16000
$ws.Cells.Item(121, 15).Value = 17000
$ws.Cells.Item(121, 16).Value = 16556
$ws.Cells.Item(121, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(121, 18).Value = "Ecuador"
$ws.Cells.Item(121, 19).Value = 1183
$ws.Cells.Item(121, 20).Value = 14
